# Weekly refresh of the "Hortaliza, Vega Monumental Concepción - Alcachofa"
# price-log sheet: the whole data block (rows 2-25) is replaced with the
# latest extract, which re-orders the previously existing rows and adds one
# new row (row 25) for the week's new entry. Columns: A Mercado ID,
# B Mercado, C Región, D Fecha, E Codreg, F Categoría ID, G Categoría,
# H Variedad, I Calidad, J Volumen, K Precio mínimo, L Precio máximo,
# M Precio promedio ponderado, N Unidad de comercialización, O Origen,
# P Precio $/Kg, Q Kg o Unidades, R Clasificación.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,18

$data[0,0] = 11
$data[0,1] = 'Vega Monumental Concepción'
$data[0,2] = 'Bíobío'
$data[0,3] = 44335
$data[0,4] = 8
$data[0,5] = 100112013
$data[0,6] = 'Alcachofa'
$data[0,7] = 'Española'
$data[0,8] = 'Primera'
$data[0,9] = 100
$data[0,10] = 17000
$data[0,11] = 18000
$data[0,12] = 17500
$data[0,13] = '$/caja 30 unidades'
$data[0,14] = 'Provincia de Limarí'
$data[0,15] = 583
$data[0,16] = 30
$data[0,17] = 'Hortaliza'

$data[1,0] = 11
$data[1,1] = 'Vega Monumental Concepción'
$data[1,2] = 'Bíobío'
$data[1,3] = 44428
$data[1,4] = 8
$data[1,5] = 100112013
$data[1,6] = 'Alcachofa'
$data[1,7] = 'Española'
$data[1,8] = 'Primera'
$data[1,9] = 100
$data[1,10] = 14000
$data[1,11] = 15000
$data[1,12] = 14500
$data[1,13] = '$/caja 30 unidades'
$data[1,14] = 'Provincia de Limarí'
$data[1,15] = 483
$data[1,16] = 30
$data[1,17] = 'Hortaliza'

$data[2,0] = 11
$data[2,1] = 'Vega Monumental Concepción'
$data[2,2] = 'Bíobío'
$data[2,3] = 44454
$data[2,4] = 8
$data[2,5] = 100112013
$data[2,6] = 'Alcachofa'
$data[2,7] = 'Madrigal'
$data[2,8] = 'Primera'
$data[2,9] = 100
$data[2,10] = 13000
$data[2,11] = 14000
$data[2,12] = 13500
$data[2,13] = '$/caja 40 unidades'
$data[2,14] = 'Provincia del Elquí'
$data[2,15] = 338
$data[2,16] = 40
$data[2,17] = 'Hortaliza'

$data[3,0] = 11
$data[3,1] = 'Vega Monumental Concepción'
$data[3,2] = 'Bíobío'
$data[3,3] = 44421
$data[3,4] = 8
$data[3,5] = 100112013
$data[3,6] = 'Alcachofa'
$data[3,7] = 'Española'
$data[3,8] = 'Primera'
$data[3,9] = 100
$data[3,10] = 14000
$data[3,11] = 15000
$data[3,12] = 14500
$data[3,13] = '$/caja 30 unidades'
$data[3,14] = 'Provincia de Limarí'
$data[3,15] = 483
$data[3,16] = 30
$data[3,17] = 'Hortaliza'

$data[4,0] = 11
$data[4,1] = 'Vega Monumental Concepción'
$data[4,2] = 'Bíobío'
$data[4,3] = 44364
$data[4,4] = 8
$data[4,5] = 100112013
$data[4,6] = 'Alcachofa'
$data[4,7] = 'Argentina(o)'
$data[4,8] = 'Primera'
$data[4,9] = 100
$data[4,10] = 19000
$data[4,11] = 20000
$data[4,12] = 19500
$data[4,13] = '$/caja 50 unidades'
$data[4,14] = 'Provincia de Limarí'
$data[4,15] = 390
$data[4,16] = 50
$data[4,17] = 'Hortaliza'

$data[5,0] = 11
$data[5,1] = 'Vega Monumental Concepción'
$data[5,2] = 'Bíobío'
$data[5,3] = 44364
$data[5,4] = 8
$data[5,5] = 100112013
$data[5,6] = 'Alcachofa'
$data[5,7] = 'Española'
$data[5,8] = 'Primera'
$data[5,9] = 100
$data[5,10] = 19000
$data[5,11] = 20000
$data[5,12] = 19500
$data[5,13] = '$/caja 30 unidades'
$data[5,14] = 'Provincia de Limarí'
$data[5,15] = 650
$data[5,16] = 30
$data[5,17] = 'Hortaliza'

$data[6,0] = 11
$data[6,1] = 'Vega Monumental Concepción'
$data[6,2] = 'Bíobío'
$data[6,3] = 44435
$data[6,4] = 8
$data[6,5] = 100112013
$data[6,6] = 'Alcachofa'
$data[6,7] = 'Argentina(o)'
$data[6,8] = 'Primera'
$data[6,9] = 100
$data[6,10] = 14000
$data[6,11] = 15000
$data[6,12] = 14500
$data[6,13] = '$/caja 50 unidades'
$data[6,14] = 'Provincia de Limarí'
$data[6,15] = 290
$data[6,16] = 50
$data[6,17] = 'Hortaliza'

$data[7,0] = 11
$data[7,1] = 'Vega Monumental Concepción'
$data[7,2] = 'Bíobío'
$data[7,3] = 44385
$data[7,4] = 8
$data[7,5] = 100112013
$data[7,6] = 'Alcachofa'
$data[7,7] = 'Española'
$data[7,8] = 'Primera'
$data[7,9] = 100
$data[7,10] = 17000
$data[7,11] = 18000
$data[7,12] = 17500
$data[7,13] = '$/caja 30 unidades'
$data[7,14] = 'Provincia de Limarí'
$data[7,15] = 583
$data[7,16] = 30
$data[7,17] = 'Hortaliza'

$data[8,0] = 11
$data[8,1] = 'Vega Monumental Concepción'
$data[8,2] = 'Bíobío'
$data[8,3] = 44342
$data[8,4] = 8
$data[8,5] = 100112013
$data[8,6] = 'Alcachofa'
$data[8,7] = 'Española'
$data[8,8] = 'Primera'
$data[8,9] = 100
$data[8,10] = 17000
$data[8,11] = 18000
$data[8,12] = 17500
$data[8,13] = '$/caja 30 unidades'
$data[8,14] = 'Provincia de Limarí'
$data[8,15] = 583
$data[8,16] = 30
$data[8,17] = 'Hortaliza'

$data[9,0] = 11
$data[9,1] = 'Vega Monumental Concepción'
$data[9,2] = 'Bíobío'
$data[9,3] = 44342
$data[9,4] = 8
$data[9,5] = 100112013
$data[9,6] = 'Alcachofa'
$data[9,7] = 'Madrigal'
$data[9,8] = 'Primera'
$data[9,9] = 100
$data[9,10] = 15000
$data[9,11] = 16000
$data[9,12] = 15500
$data[9,13] = '$/caja 40 unidades'
$data[9,14] = 'Provincia de Limarí'
$data[9,15] = 388
$data[9,16] = 40
$data[9,17] = 'Hortaliza'

$data[10,0] = 11
$data[10,1] = 'Vega Monumental Concepción'
$data[10,2] = 'Bíobío'
$data[10,3] = 44399
$data[10,4] = 8
$data[10,5] = 100112013
$data[10,6] = 'Alcachofa'
$data[10,7] = 'Española'
$data[10,8] = 'Primera'
$data[10,9] = 100
$data[10,10] = 14000
$data[10,11] = 15000
$data[10,12] = 14500
$data[10,13] = '$/caja 30 unidades'
$data[10,14] = 'Provincia de Limarí'
$data[10,15] = 483
$data[10,16] = 30
$data[10,17] = 'Hortaliza'

$data[11,0] = 11
$data[11,1] = 'Vega Monumental Concepción'
$data[11,2] = 'Bíobío'
$data[11,3] = 44383
$data[11,4] = 8
$data[11,5] = 100112013
$data[11,6] = 'Alcachofa'
$data[11,7] = 'Argentina(o)'
$data[11,8] = 'Primera'
$data[11,9] = 50
$data[11,10] = 17000
$data[11,11] = 18000
$data[11,12] = 17400
$data[11,13] = '$/caja 50 unidades'
$data[11,14] = 'Provincia de Limarí'
$data[11,15] = 348
$data[11,16] = 50
$data[11,17] = 'Hortaliza'

$data[12,0] = 11
$data[12,1] = 'Vega Monumental Concepción'
$data[12,2] = 'Bíobío'
$data[12,3] = 44441
$data[12,4] = 8
$data[12,5] = 100112013
$data[12,6] = 'Alcachofa'
$data[12,7] = 'Española'
$data[12,8] = 'Primera'
$data[12,9] = 100
$data[12,10] = 13000
$data[12,11] = 14000
$data[12,12] = 13500
$data[12,13] = '$/caja 30 unidades'
$data[12,14] = 'Provincia de Limarí'
$data[12,15] = 450
$data[12,16] = 30
$data[12,17] = 'Hortaliza'

$data[13,0] = 11
$data[13,1] = 'Vega Monumental Concepción'
$data[13,2] = 'Bíobío'
$data[13,3] = 44426
$data[13,4] = 8
$data[13,5] = 100112013
$data[13,6] = 'Alcachofa'
$data[13,7] = 'Madrigal'
$data[13,8] = 'Primera'
$data[13,9] = 50
$data[13,10] = 12000
$data[13,11] = 13000
$data[13,12] = 12600
$data[13,13] = '$/caja 40 unidades'
$data[13,14] = 'Provincia de Limarí'
$data[13,15] = 315
$data[13,16] = 40
$data[13,17] = 'Hortaliza'

$data[14,0] = 11
$data[14,1] = 'Vega Monumental Concepción'
$data[14,2] = 'Bíobío'
$data[14,3] = 44442
$data[14,4] = 8
$data[14,5] = 100112013
$data[14,6] = 'Alcachofa'
$data[14,7] = 'Española'
$data[14,8] = 'Primera'
$data[14,9] = 100
$data[14,10] = 14500
$data[14,11] = 15000
$data[14,12] = 14750
$data[14,13] = '$/caja 30 unidades'
$data[14,14] = 'Provincia de Limarí'
$data[14,15] = 492
$data[14,16] = 30
$data[14,17] = 'Hortaliza'

$data[15,0] = 11
$data[15,1] = 'Vega Monumental Concepción'
$data[15,2] = 'Bíobío'
$data[15,3] = 44420
$data[15,4] = 8
$data[15,5] = 100112013
$data[15,6] = 'Alcachofa'
$data[15,7] = 'Española'
$data[15,8] = 'Primera'
$data[15,9] = 100
$data[15,10] = 14000
$data[15,11] = 15000
$data[15,12] = 14500
$data[15,13] = '$/caja 30 unidades'
$data[15,14] = 'Provincia de Limarí'
$data[15,15] = 483
$data[15,16] = 30
$data[15,17] = 'Hortaliza'

$data[16,0] = 11
$data[16,1] = 'Vega Monumental Concepción'
$data[16,2] = 'Bíobío'
$data[16,3] = 44350
$data[16,4] = 8
$data[16,5] = 100112013
$data[16,6] = 'Alcachofa'
$data[16,7] = 'Argentina(o)'
$data[16,8] = 'Primera'
$data[16,9] = 50
$data[16,10] = 15000
$data[16,11] = 16000
$data[16,12] = 15600
$data[16,13] = '$/caja 50 unidades'
$data[16,14] = 'Provincia de Limarí'
$data[16,15] = 312
$data[16,16] = 50
$data[16,17] = 'Hortaliza'

$data[17,0] = 11
$data[17,1] = 'Vega Monumental Concepción'
$data[17,2] = 'Bíobío'
$data[17,3] = 44350
$data[17,4] = 8
$data[17,5] = 100112013
$data[17,6] = 'Alcachofa'
$data[17,7] = 'Española'
$data[17,8] = 'Primera'
$data[17,9] = 40
$data[17,10] = 17000
$data[17,11] = 18000
$data[17,12] = 17500
$data[17,13] = '$/caja 30 unidades'
$data[17,14] = 'Provincia de Limarí'
$data[17,15] = 583
$data[17,16] = 30
$data[17,17] = 'Hortaliza'

$data[18,0] = 11
$data[18,1] = 'Vega Monumental Concepción'
$data[18,2] = 'Bíobío'
$data[18,3] = 44447
$data[18,4] = 8
$data[18,5] = 100112013
$data[18,6] = 'Alcachofa'
$data[18,7] = 'Española'
$data[18,8] = 'Primera'
$data[18,9] = 100
$data[18,10] = 14000
$data[18,11] = 15000
$data[18,12] = 14500
$data[18,13] = '$/caja 30 unidades'
$data[18,14] = 'Provincia de Limarí'
$data[18,15] = 483
$data[18,16] = 30
$data[18,17] = 'Hortaliza'

$data[19,0] = 11
$data[19,1] = 'Vega Monumental Concepción'
$data[19,2] = 'Bíobío'
$data[19,3] = 44433
$data[19,4] = 8
$data[19,5] = 100112013
$data[19,6] = 'Alcachofa'
$data[19,7] = 'Argentina(o)'
$data[19,8] = 'Primera'
$data[19,9] = 100
$data[19,10] = 14000
$data[19,11] = 15000
$data[19,12] = 14500
$data[19,13] = '$/caja 50 unidades'
$data[19,14] = 'Provincia de Limarí'
$data[19,15] = 290
$data[19,16] = 50
$data[19,17] = 'Hortaliza'

$data[20,0] = 11
$data[20,1] = 'Vega Monumental Concepción'
$data[20,2] = 'Bíobío'
$data[20,3] = 44397
$data[20,4] = 8
$data[20,5] = 100112013
$data[20,6] = 'Alcachofa'
$data[20,7] = 'Española'
$data[20,8] = 'Primera'
$data[20,9] = 100
$data[20,10] = 14000
$data[20,11] = 15000
$data[20,12] = 14500
$data[20,13] = '$/caja 30 unidades'
$data[20,14] = 'Provincia de Limarí'
$data[20,15] = 483
$data[20,16] = 30
$data[20,17] = 'Hortaliza'

$data[21,0] = 11
$data[21,1] = 'Vega Monumental Concepción'
$data[21,2] = 'Bíobío'
$data[21,3] = 44376
$data[21,4] = 8
$data[21,5] = 100112013
$data[21,6] = 'Alcachofa'
$data[21,7] = 'Española'
$data[21,8] = 'Primera'
$data[21,9] = 100
$data[21,10] = 19000
$data[21,11] = 20000
$data[21,12] = 19500
$data[21,13] = '$/caja 30 unidades'
$data[21,14] = 'Provincia de Limarí'
$data[21,15] = 650
$data[21,16] = 30
$data[21,17] = 'Hortaliza'

$data[22,0] = 11
$data[22,1] = 'Vega Monumental Concepción'
$data[22,2] = 'Bíobío'
$data[22,3] = 44358
$data[22,4] = 8
$data[22,5] = 100112013
$data[22,6] = 'Alcachofa'
$data[22,7] = 'Argentina(o)'
$data[22,8] = 'Primera'
$data[22,9] = 100
$data[22,10] = 18000
$data[22,11] = 20000
$data[22,12] = 19000
$data[22,13] = '$/caja 50 unidades'
$data[22,14] = 'Provincia de Limarí'
$data[22,15] = 380
$data[22,16] = 50
$data[22,17] = 'Hortaliza'

$data[23,0] = 11
$data[23,1] = 'Vega Monumental Concepción'
$data[23,2] = 'Bíobío'
$data[23,3] = 44358
$data[23,4] = 8
$data[23,5] = 100112013
$data[23,6] = 'Alcachofa'
$data[23,7] = 'Española'
$data[23,8] = 'Primera'
$data[23,9] = 100
$data[23,10] = 18000
$data[23,11] = 20000
$data[23,12] = 19000
$data[23,13] = '$/caja 30 unidades'
$data[23,14] = 'Provincia de Limarí'
$data[23,15] = 633
$data[23,16] = 30
$data[23,17] = 'Hortaliza'

$ws.Range("A2:R25").Value = $data

# Row 25 is a brand-new row; give its date cell (column D) the same
# date number format used by the rest of column D (style index 2).
$ws.Range("D25").NumberFormat = $ws.Range("D24").NumberFormat
